# Battery_Data.xlsx update
# - Update the 4 result values on the "Battery_Data" sheet (sheet1)
# - Update the 4 existing yearly values and append 16 new yearly rows
#   (rows 6-21, years 5-20) on the "Yearly BRC" sheet (sheet2)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Battery_Data
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Battery_Data")

$ws1.Range("B2").Value = 34.3685988048
$ws1.Range("B3").Value = 18.559043354592
$ws1.Range("B4").Value = 0.3711808670918401
$ws1.Range("B5").Value = 6.37166537488

# ---------------------------------------------------------------------
# Sheet 2: Yearly BRC
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Yearly BRC")

# Update the existing 4 rows (years 1-4)
$ws2.Range("B2").Value = 0.7982541322620772
$ws2.Range("B3").Value = 0.81737474436244
$ws2.Range("B4").Value = 0.8383750009946094
$ws2.Range("B5").Value = 0.8620796584972996

# Extend the formatting of column A down to row 21 (same style as the
# existing label cells) before filling in the new labels/values.
$ws2.Range("A5").Copy($ws2.Range("A6:A21"))

# New rows for years 5-20
$newValues = @{
    6  = 0.8897504170499579
    7  = 0.9225060407143506
    8  = 0.9619115834189335
    9  = 1.0095103835637
    10 = 1.06531849117498
    11 = 1.129309976919592
    12 = 1.201802342693103
    13 = 1.282820757874964
    14 = 1.372014508240509
    15 = 1.468652028395985
    16 = 1.572213123728664
    17 = 1.682267726638863
    18 = 1.798366570241442
    19 = 1.919892149598082
    20 = 2.045897355990393
    21 = 2.175136961007468
}

for ($row = 6; $row -le 21; $row++) {
    $year = $row - 1
    $ws2.Range("A$row").Value = "Battery Replacement Cost at y = $year"
    $ws2.Range("B$row").Value = $newValues[$row]
}
